# Daily attendance processing - 2025-11-25 21:48:09
# For every cell in column "G" (Recorded By) whose value starts with the
# literal prefix "System, ", move the "System" token from the front of the
# comma-separated list to the end (e.g. "System, admin@admin.com" becomes
# "admin@admin.com, System"). Values that don't start with that exact
# prefix (e.g. a bare "System", or values where "System" already trails,
# or a lower-case "system, ...") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$prefix = "System, "

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = $cell.Text

    if ($text -and $text.StartsWith($prefix)) {
        $rest = $text.Substring($prefix.Length)
        $newValue = $rest + ", System"
        $cell.Value = $newValue
    }
}
